$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "tunnit yht." totals row (row 60); it will be re-created
#    further down the sheet (row 75) once the new entries are in place.
# ---------------------------------------------------------------------------
$ws.Rows.Item(60).Delete()

# ---------------------------------------------------------------------------
# 2. New diary entries (rows 59-64). Styles are picked up from neighbouring
#    cells via Copy so the shared cellXfs indices line up exactly like they
#    would if typed by hand in Excel.
# ---------------------------------------------------------------------------
$ws.Range("A57").Copy($ws.Range("A59"))
$ws.Range("A57").Copy($ws.Range("A60"))
$ws.Range("B57:D57").Copy($ws.Range("B59:D59"))
$ws.Range("B57:D57").Copy($ws.Range("B60:D60"))
$ws.Range("B58:D58").Copy($ws.Range("B61:D61"))
$ws.Range("B58:D58").Copy($ws.Range("B62:D62"))
$ws.Range("B58:D58").Copy($ws.Range("B63:D63"))
$ws.Range("B58:D58").Copy($ws.Range("B64:D64"))

$ws.Cells.Item(59,1).Value = 44565
$ws.Cells.Item(59,2).Value = 1
$ws.Cells.Item(59,3).Value = "login form luotu, router lisätty jotta navigointia voidaan käyttää, analyzer container datan hakuun luotu"
$ws.Cells.Item(59,4).Value = "client"

$ws.Cells.Item(60,1).Value = 44567
$ws.Cells.Item(60,2).Value = 1
$ws.Cells.Item(60,3).Value = "redux perusasetukset tehty, ekan analyzerReducer rakentamista"
$ws.Cells.Item(60,4).Value = "client"

$ws.Cells.Item(61,2).Value = 2
$ws.Cells.Item(61,3).Value = "error reducer tehty ja koodin refaktorointia ja testausta"
$ws.Cells.Item(61,4).Value = "client"

$ws.Cells.Item(62,2).Value = 1
$ws.Cells.Item(62,3).Value = "login form css ja statet muokattu"
$ws.Cells.Item(62,4).Value = "client"

$ws.Cells.Item(63,2).Value = 1
$ws.Cells.Item(63,3).Value = "oma uudelleenkäytettävä komponentti formien submit/cancel napeille, css luokkien uudelleen nimeämistä"
$ws.Cells.Item(63,4).Value = "client"

$ws.Cells.Item(64,2).Value = 1
$ws.Cells.Item(64,3).Value = "loginservice luotu, ei testattu, userReducer ja lisäys storeen"
$ws.Cells.Item(64,4).Value = "client"

# ---------------------------------------------------------------------------
# 3. Re-create the totals row, now at row 75, summing the extended range.
#    Values/formula are written first and the look (number format / align-
#    ment) is applied afterwards via a formats-only paste - doing the style
#    paste before the formula corrupts the live formula dependency graph.
#    C1 carries the same style as the old A60 ("tunnit yht." label: centred
#    shared-string cell); B57 carries the plain numeric style used by B60.
# ---------------------------------------------------------------------------
$ws.Cells.Item(75,1).Value = "tunnit yht."
$ws.Cells.Item(75,2).Formula = "=SUM(B2:B64)"
$ws.Range("C1").Copy()
$ws.Range("A75").PasteSpecial(-4122)
$ws.Range("B57").Copy()
$ws.Range("B75").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Sheet view bookkeeping: selection mirrors the author's final cursor
#    position in the workbook.
# ---------------------------------------------------------------------------
$ws.Range("D64").Select()
